$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all data rows (2-32)
# from 45207 to 45208.
for ($row = 2; $row -le 32; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}

# Update the HYPERLINK formulas in columns S, T, V, W, X, Y for rows 2-5,
# replacing the folder name "Logging_ALVDALEN" with "Logging_2039".
$hyperlinkColumns = @("S", "T", "V", "W", "X", "Y")
for ($row = 2; $row -le 5; $row++) {
    foreach ($col in $hyperlinkColumns) {
        $cell = $ws.Range("$col$row")
        $formula = $cell.Formula
        if ($formula -and $formula.Contains("Logging_ALVDALEN")) {
            $cell.Formula = $formula.Replace("Logging_ALVDALEN", "Logging_2039")
        }
    }
}
